$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9642150815792832
$ws.Range("C2").Value = 0.1235905131834478
$ws.Range("D2").Value = 0.5483649432170807
$ws.Range("E2").Value = 0.1864854176659154
$ws.Range("G2").Value = 0.002592011712660325
$ws.Range("J2").Value = 0.07734300115510706
$ws.Range("K2").Value = 0.4534032876824483
$ws.Range("L2").Value = 0.3725307360031636
$ws.Range("M2").Value = 0.296807011549717
$ws.Range("O2").Value = 9.087362009319634

$ws.Range("B3").Value = 0.9385773748042823
$ws.Range("C3").Value = 0.1227812568374667
$ws.Range("D3").Value = 0.5469674629464976
$ws.Range("E3").Value = 0.1869680779558625
$ws.Range("G3").Value = 0.002595137568200568
$ws.Range("J3").Value = 0.07743948547543678
$ws.Range("K3").Value = 0.4292702986610237
$ws.Range("L3").Value = 0.3712930136328865
$ws.Range("M3").Value = 0.2923225972467556
$ws.Range("O3").Value = 9.108522568869432

$ws.Range("B4").Value = 0.9233060635167476
$ws.Range("C4").Value = 0.1222729659786204
$ws.Range("D4").Value = 0.546318507956002
$ws.Range("E4").Value = 0.187315391219439
$ws.Range("G4").Value = 0.002597160417498293
$ws.Range("J4").Value = 0.07750173979297115
$ws.Range("K4").Value = 0.4146744601719092
$ws.Range("L4").Value = 0.3706713037927472
$ws.Range("M4").Value = 0.2897035274509356
$ws.Range("O4").Value = 9.124641079459025

$ws.Range("B5").Value = 0.9172016949162582
$ws.Range("C5").Value = 0.1220629682301499
$ws.Range("D5").Value = 0.5461067475730061
$ws.Range("E5").Value = 0.1874697652377222
$ws.Range("G5").Value = 0.002598010864801457
$ws.Range("J5").Value = 0.07752786834406411
$ws.Range("K5").Value = 0.408782711452389
$ws.Range("L5").Value = 0.3704528115507273
$ws.Range("M5").Value = 0.2886701506551965
$ws.Range("O5").Value = 9.131996202060265

$ws.Range("B6").Value = 0.9161952613106905
$ws.Range("C6").Value = 0.1220279252984113
$ws.Range("D6").Value = 0.5460747708522717
$ws.Range("E6").Value = 0.1874961752959958
$ws.Range("G6").Value = 0.002598153660930955
$ws.Range("J6").Value = 0.07753225289181787
$ws.Range("K6").Value = 0.4078077948295089
$ws.Range("L6").Value = 0.3704186392784266
$ws.Range("M6").Value = 0.2885006113163868
$ws.Range("O6").Value = 9.133265049317032

$ws.Range("B7").Value = 0.9232232560682689
$ws.Range("C7").Value = 0.1222701454728217
$ws.Range("D7").Value = 0.5463154385745668
$ws.Range("E7").Value = 0.1873174211321622
$ws.Range("G7").Value = 0.002597171781042334
$ws.Range("J7").Value = 0.07750208909456102
$ws.Range("K7").Value = 0.4145947740472593
$ws.Range("L7").Value = 0.3706682158596308
$ws.Range("M7").Value = 0.2896894534817775
$ws.Range("O7").Value = 9.124737087016399

$ws.Range("B8").Value = 0.9552778948570051
$ws.Range("C8").Value = 0.1233138498620221
$ws.Range("D8").Value = 0.5478397588989168
$ws.Range("E8").Value = 0.1866412788253449
$ws.Range("G8").Value = 0.002593068060996049
$ws.Range("J8").Value = 0.07737564482842529
$ws.Range("K8").Value = 0.4450363886205793
$ws.Range("L8").Value = 0.3720753291862238
$ws.Range("M8").Value = 0.2952329678040897
$ws.Range("O8").Value = 9.094009799023439

$ws.Range("B9").Value = 1.02184940205774
$ws.Range("C9").Value = 0.1252700286891866
$ws.Range("D9").Value = 0.5524837845300112
$ws.Range("E9").Value = 0.1857185579671903
$ws.Range("G9").Value = 0.00258583879918753
$ws.Range("J9").Value = 0.0771515043955695
$ws.Range("K9").Value = 0.5064803196257515
$ws.Range("L9").Value = 0.3759280458096441
$ws.Range("M9").Value = 0.3071654181459316
$ws.Range("O9").Value = 9.058530975147846

$ws.Range("B10").Value = 1.073002872521869
$ws.Range("C10").Value = 0.1266521165130854
$ws.Range("D10").Value = 0.556899789014949
$ws.Range("E10").Value = 0.1852849962749801
$ws.Range("G10").Value = 0.002581021179767422
$ws.Range("J10").Value = 0.07700122434405188
$ws.Range("K10").Value = 0.552676902103201
$ws.Range("L10").Value = 0.3794210212791427
$ws.Range("M10").Value = 0.3165745125208588
$ws.Range("O10").Value = 9.047543650333893

$ws.Range("B11").Value = 1.096757272839852
$ws.Range("C11").Value = 0.1272689189286993
$ws.Range("D11").Value = 0.5591258144697235
$ws.Range("E11").Value = 0.1851405224463285
$ws.Range("G11").Value = 0.002578935656736963
$ws.Range("J11").Value = 0.07693595783102891
$ws.Range("K11").Value = 0.5739195598364404
$ws.Range("L11").Value = 0.3811530667619678
$ws.Range("M11").Value = 0.320993483687424
$ws.Range("O11").Value = 9.045814597681272

$ws.Range("B12").Value = 1.105821658185562
$ws.Range("C12").Value = 0.1275007717481529
$ws.Range("D12").Value = 0.5599998825076824
$ws.Range("E12").Value = 0.1850933745512187
$ws.Range("G12").Value = 0.002578161090319187
$ws.Range("J12").Value = 0.07691168647630864
$ws.Range("K12").Value = 0.5819960358697074
$ws.Range("L12").Value = 0.3818294386897918
$ws.Range("M12").Value = 0.3226866766057626
$ws.Range("O12").Value = 9.0456294310776

$ws.Range("B13").Value = 1.103866415702157
$ws.Range("C13").Value = 0.1274509145044007
$ws.Range("D13").Value = 0.5598102539277079
$ws.Range("E13").Value = 0.1851031928038793
$ws.Range("G13").Value = 0.002578327233358297
$ws.Range("J13").Value = 0.07691689404076874
$ws.Range("K13").Value = 0.5802551899311084
$ws.Range("L13").Value = 0.3816828602935374
$ws.Range("M13").Value = 0.3223211374921462
$ws.Range("O13").Value = 9.04564843304064

$ws.Range("B14").Value = 1.097501623367805
$ws.Range("C14").Value = 0.1272880280271096
$ws.Range("D14").Value = 0.5591971014843722
$ws.Range("E14").Value = 0.1851364921653733
$ws.Range("G14").Value = 0.002578871628857318
$ws.Range("J14").Value = 0.07693395212998055
$ws.Range("K14").Value = 0.5745833703916219
$ws.Range("L14").Value = 0.3812083022733077
$ws.Range("M14").Value = 0.3211323870866494
$ws.Range("O14").Value = 9.045789955577618

$ws.Range("B15").Value = 1.093611988068659
$ws.Range("C15").Value = 0.1271880316718281
$ws.Range("D15").Value = 0.5588255775388262
$ws.Range("E15").Value = 0.1851578729390013
$ws.Range("G15").Value = 0.002579207061304046
$ws.Range("J15").Value = 0.07694445843131437
$ws.Range("K15").Value = 0.5711134176856376
$ws.Range("L15").Value = 0.3809202867296904
$ws.Range("M15").Value = 0.3204068211861468
$ws.Range("O15").Value = 9.045937781373709

$ws.Range("B16").Value = 1.071460167410891
$ws.Range("C16").Value = 0.1266115669530876
$ws.Range("D16").Value = 0.5567586737138157
$ws.Range("E16").Value = 0.1852954974185685
$ws.Range("G16").Value = 0.00258115960094443
$ws.Range("J16").Value = 0.07700555183599267
$ws.Range("K16").Value = 0.5512931913167733
$ws.Range("L16").Value = 0.3793106993429376
$ws.Range("M16").Value = 0.3162885032630669
$ws.Range("O16").Value = 9.047722387729635

$ws.Range("B17").Value = 1.057994424259789
$ws.Range("C17").Value = 0.1262548700770978
$ws.Range("D17").Value = 0.5555462393815844
$ws.Range("E17").Value = 0.185393419735977
$ws.Range("G17").Value = 0.002582384527988424
$ws.Range("J17").Value = 0.07704382261389053
$ws.Range("K17").Value = 0.5391921430339721
$ws.Range("L17").Value = 0.37835985006069
$ws.Range("M17").Value = 0.3137974912706483
$ws.Range("O17").Value = 9.049654126564405

$ws.Range("B18").Value = 1.050294914578075
$ws.Range("C18").Value = 0.1260485860455489
$ws.Range("D18").Value = 0.5548693260131188
$ws.Range("E18").Value = 0.1854547091522711
$ws.Range("G18").Value = 0.002583099058672168
$ws.Range("J18").Value = 0.07706612651980826
$ws.Range("K18").Value = 0.5322533859290957
$ws.Range("L18").Value = 0.3778264187985627
$ws.Range("M18").Value = 0.3123777920909276
$ws.Range("O18").Value = 9.051072948460671

$ws.Range("B19").Value = 1.047695846657717
$ws.Range("C19").Value = 0.1259785493136789
$ws.Range("D19").Value = 0.5546436498061809
$ws.Range("E19").Value = 0.1854763145608054
$ws.Range("G19").Value = 0.002583342703561139
$ws.Range("J19").Value = 0.07707372836620063
$ws.Range("K19").Value = 0.5299077371332999
$ws.Range("L19").Value = 0.3776481247297454
$ws.Range("M19").Value = 0.3118993540719188
$ws.Range("O19").Value = 9.051606206878603

$ws.Range("B20").Value = 1.0594231571107
$ws.Range("C20").Value = 0.1262929571561244
$ws.Range("D20").Value = 0.5556731896984388
$ws.Range("E20").Value = 0.185382481844222
$ws.Range("G20").Value = 0.00258225309963995
$ws.Range("J20").Value = 0.07703971845989699
$ws.Range("K20").Value = 0.5404781035252597
$ws.Range("L20").Value = 0.3784596759312109
$ws.Range("M20").Value = 0.3140613124703791
$ws.Range("O20").Value = 9.049416643745502

$ws.Range("B21").Value = 1.099369245927875
$ws.Range("C21").Value = 0.1273359183550511
$ws.Range("D21").Value = 0.5593763554716418
$ws.Range("E21").Value = 0.1851265063252612
$ws.Range("G21").Value = 0.002578711315416962
$ws.Range("J21").Value = 0.07692892972781395
$ws.Range("K21").Value = 0.5762484465764146
$ws.Range("L21").Value = 0.3813471362547034
$ws.Range("M21").Value = 0.3214810148548821
$ws.Range("O21").Value = 9.045735646537651

$ws.Range("B22").Value = 1.125878893237683
$ws.Range("C22").Value = 0.1280075450996421
$ws.Range("D22").Value = 0.561977932961824
$ws.Range("E22").Value = 0.1850032725238933
$ws.Range("G22").Value = 0.002576484978746098
$ws.Range("J22").Value = 0.07685910811489904
$ws.Range("K22").Value = 0.5998147730867345
$ws.Range("L22").Value = 0.3833536128851023
$ws.Range("M22").Value = 0.3264457140017285
$ws.Range("O22").Value = 9.046066786659424

$ws.Range("B23").Value = 1.111693541849149
$ws.Range("C23").Value = 0.1276500021313893
$ws.Range("D23").Value = 0.5605728644097212
$ws.Range("E23").Value = 0.1850650213617691
$ws.Range("G23").Value = 0.002577665149329974
$ws.Range("J23").Value = 0.07689613721033517
$ws.Range("K23").Value = 0.5872198750909661
$ws.Range("L23").Value = 0.3822718272024872
$ws.Range("M23").Value = 0.3237854337890269
$ws.Range("O23").Value = 9.045639793455052

$ws.Range("B24").Value = 1.058777095731472
$ws.Range("C24").Value = 0.1262757417695326
$ws.Range("D24").Value = 0.555615732746233
$ws.Range("E24").Value = 0.1853874113120213
$ws.Range("G24").Value = 0.002582312486344181
$ws.Range("J24").Value = 0.07704157300878967
$ws.Range("K24").Value = 0.5398966637667115
$ws.Range("L24").Value = 0.3784145034553319
$ws.Range("M24").Value = 0.3139420002198321
$ws.Range("O24").Value = 9.049523049512715

$ws.Range("B25").Value = 1.003444661546382
$ws.Range("C25").Value = 0.1247505149092873
$ws.Range("D25").Value = 0.5510507876538782
$ws.Range("E25").Value = 0.1859251627373357
$ws.Range("G25").Value = 0.002587707442600854
$ws.Range("J25").Value = 0.07720960342261973
$ws.Range("K25").Value = 0.4896722105652316
$ws.Range("L25").Value = 0.3747691630352818
$ws.Range("M25").Value = 0.3038242253445524
$ws.Range("O25").Value = 9.065479051996135
